$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 5 (pushes old row 5 "LED 1206" data and the
# trailing blank/styled rows down by two, and shifts the totals row too).
[void]$ws.Rows("5:6").Insert()

# Row 2 ("N/A" / PCB line): the unit price is no longer filled in yet,
# so clear it; the shared formula recalculates to 0.
$ws.Range("D2").ClearContents()

# Row 4 (existing Harting header, 1106765): more specific description.
$ws.Range("F4").Value = "Hartin Connector Chassis"

# Row 5 (NEW): Harting connector cable part.
$ws.Range("A5").Value = 1106731
$ws.Range("B5").Value = "Farnell"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2.5499999999999998
$ws.Range("E5").Formula = "=C5*D5"
$ws.Range("F5").Value = "Hartin Connector Cable"
[void]$ws.Hyperlinks.Add($ws.Range("A5"), "https://be.farnell.com/harting/09-18-526-6928/header-straight-s-latch-26way/dp/1106731", "", "", "1106731")
# Adding the hyperlink resets the cell style; restore the plain "hyperlink" xf
# that the rest of the part-number column already uses.
$ws.Range("A5").Style = $ws.Range("A4").Style

# Row 6 (NEW): ribbon cable, amount only (price still TBD).
$ws.Range("B6").Value = "Farnell"
$ws.Range("C6").Value = 1
$ws.Range("F6").Value = "25-26 ribbon cable."

# Row 7 (previously row 5, LED 1206 line) keeps its numbers but the
# description text moved position in the shared-string table.
$ws.Range("F7").Value = "LED 1206"

# Column F needs to be a bit wider to fit the new, longer descriptions.
$ws.Columns("F").ColumnWidth = 22.43

# Leave the selection where the author left off.
[void]$ws.Range("E8").Select()
